# Fix "before present" projection times/selections across Station1-Station5 sheets.
$wb = $excel.ActiveWorkbook

# --- Station1 (sheet1) ---
$ws1 = $wb.Worksheets.Item("Station1")
$ws1.Range("B2").Value = 0.375
$ws1.Activate()
$ws1.Range("E14").Select()

# --- Station2 (sheet2) ---
$ws2 = $wb.Worksheets.Item("Station2")
$ws2.Range("B2").Value = 0.375
$ws2.Activate()
$ws2.Range("B3").Select()

# --- Station3 (sheet3) ---
$ws3 = $wb.Worksheets.Item("Station3")
$ws3.Range("B2").Value = 0.375
$ws3.Range("B3").Value = 0.41666666666666669
$ws3.Activate()
$ws3.Range("B4").Select()

# --- Station4 (sheet4) ---
$ws4 = $wb.Worksheets.Item("Station4")
$ws4.Range("B2").Value = 0.375
$ws4.Range("B3").Value = 0.41666666666666669
$ws4.Activate()
$ws4.Range("B4").Select()

# --- Station5 (sheet5) : keep as the active sheet, same as original workbook state ---
$ws5 = $wb.Worksheets.Item("Station5")
$ws5.Range("B3").Value = 0.41666666666666669
$ws5.Range("B4").Value = 0.45833333333333331
$ws5.Activate()
$ws5.Range("D15").Select()
